# Generate Report for Handback
# Updates the handoff/handback timestamps produced by a new report generation run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-01 00:52:44"
# G3 is unchanged ("2016-09-01 00:51:52")

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-01 00:52:40"
$wsZhCn.Range("K2").Value = "2016-09-01 00:52:57"
# Row 3 (fcce5ca1) is unchanged

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K2").Value = "2016-09-01 00:53:11"
$wsDeDe.Range("H3").Value = "2016-09-01 00:52:44"
# H2 and K3 are unchanged
